$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sheet view: scroll position / selection change ---
# topLeftCell moves from B1 to B10, and the active selection moves from D18 to C30.
$excel.Goto($ws.Range("B10"), $true)
[void]$ws.Range("C30").Select()

# --- Row 23: new Application/Task entries (must be written before row 18/19/22
#     edits below so the shared-string table ends up in the same order as the
#     authoritative workbook) ---
$ws.Range("C23").Value = "QMVAR"
$ws.Range("D23").Value = "Weekly Revenue Report analysis the scenario going on."

# --- Row 22: Application cell previously existed (blank, bordered) - clear its
#     formatting first so it picks up the plain column style, then fill it in ---
$ws.Range("C22").Clear()
$ws.Range("C22").Value = "QMVAR"
$ws.Range("D22").Value = "Weekly Revenue Report gather requirements and analysis the scenario going on."

# --- Row 18: Task text updated to mention "given by mohansan" ---
$ws.Range("D18").Value = "Mobile view task implemented in latest file in Qmvar given by mohansan and tested"

# --- Row 19: new Application/Task entries ---
$ws.Range("C19").Value = "QMVAR"
$ws.Range("D19").Value = "Mobile view task sent to deployment and Weekly Revenue report document study started."

# --- Row 24: new Application/Task entries ---
$ws.Range("C24").Value = "QMVAR"
$ws.Range("D24").Value = "Weekly Revenue Report db analysis and implementation going  on"

# --- Row 25: new Application/Task entries ---
$ws.Range("C25").Value = "QMVAR"
$ws.Range("D25").Value = "Weekly Revenue Report db datas checked and tested for all service centers"

# --- Row 26: new Application/Task entries ---
$ws.Range("C26").Value = "QMVAR"
$ws.Range("D26").Value = "Weekly Revenue Report Service centers calculation implementation going on"

# --- Row 29: new Application/Task entries (same task text as row 26) ---
$ws.Range("C29").Value = "QMVAR"
$ws.Range("D29").Value = "Weekly Revenue Report Service centers calculation implementation going on"
